# Crear/Actualizar Excel para pedido 69134d11b9c1d30b15fabdc3
#
# The sheet's row 11 is a recurring "Test Ringover (NO TOCAR)" QA record
# (order 2488). This adds a fresh copy of that same record as a new row 12
# (same values in every column), and tidies up the original row's blank
# "Paneles" placeholder cell in the process.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 11 (the existing "Test Ringover" record) into row 12 so the
# new test-order row carries identical data across every column.
$ws.Rows.Item(11).Copy()
$ws.Rows.Item(12).PasteSpecial()

# Row 11's "Paneles" cell (D11) was a blank placeholder (no text). The copy
# above drops truly-blank cells, so recreate that same blank placeholder on
# the new row by toggling a formatting property, which forces the cell to
# be materialized without giving it any value.
$ws.Range("D12").Font.Bold = $true
$ws.Range("D12").Font.Bold = $false

# The original row no longer needs its own blank "Paneles" placeholder cell.
$ws.Range("D11").ClearContents()
